# Update version string across the workbook for the new release:
#   "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# becomes
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: <version>"
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation sentence containing the version string
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Moonidih Coal Mine, India, M2881, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# S2:S11 on "Boundaries and methane sources" sheet hold the build_version value
for ($row = 2; $row -le 11; $row++) {
    $wsBoundaries.Cells.Item($row, 19).Value = $newVersion
}
